# Update the "Correspond Handoff Datetime" / "Correspond Handback DateTime"
# timestamps on the zh-cn and de-de report rows to reflect the newly
# generated handback report.

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E4").Value = "2016-03-11 09:33:53"
$wsZhCn.Range("H4").Value = "2016-03-11 09:34:10"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E4").Value = "2016-03-11 09:33:56"
$wsDeDe.Range("H4").Value = "2016-03-11 09:34:15"
